# prod smoke test - platform testing
# Update the Lexus MSRP staging sheet: roll several trims from MY2020 to
# MY2021 pricing, and add four new 2021 "Black Line" RC trims.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bump model-year + refresh BASE MSRP for the existing RC 300/RC 350
#     rows (2-9) and the RC F / RC F TRACK rows (53-54) ---
$yearUpdates = @(
  @{ Row = 2;  Year = 2021; Msrp = 42120 },
  @{ Row = 3;  Year = 2021; Msrp = 46590 },
  @{ Row = 4;  Year = 2021; Msrp = 44810 },
  @{ Row = 5;  Year = 2021; Msrp = 48765 },
  @{ Row = 6;  Year = 2021; Msrp = 45050 },
  @{ Row = 7;  Year = 2021; Msrp = 49520 },
  @{ Row = 8;  Year = 2021; Msrp = 47215 },
  @{ Row = 9;  Year = 2021; Msrp = 51130 },
  @{ Row = 53; Year = 2021; Msrp = 65875 },
  @{ Row = 54; Year = 2021; Msrp = 96675 }
)

foreach ($u in $yearUpdates) {
  $ws.Cells.Item($u.Row, 3).Value = $u.Year
  $ws.Cells.Item($u.Row, 4).Value = $u.Msrp
}

# --- Append four new 2021 "Black Line" RC trims below the existing data
#     (rows 95-98). Trim codes go in column A, then trim names in column
#     B, matching how the sheet was actually authored. ---
$ws.Range("A95").Value = "9203SE"
$ws.Range("A96").Value = "9207SE"
$ws.Range("A97").Value = "9213SE"
$ws.Range("A98").Value = "9217SE"

$ws.Range("B95").Value = "RC 300 F SPORT Black Line"
$ws.Range("B96").Value = "RC 300 AWD F SPORT Black Line"
$ws.Range("B97").Value = "RC 350 F SPORT Black Line"
$ws.Range("B98").Value = "RC 350 AWD F SPORT Black Line"

$ws.Range("C95:C98").Value = 2021

$ws.Range("D95").Value = 48735
$ws.Range("D96").Value = 50910
$ws.Range("D97").Value = 51665
$ws.Range("D98").Value = 53275
$ws.Range("D95:D98").NumberFormat = "$#,##0_);[Red]($#,##0)"

$ws.Range("E95:E98").Value = 1025
$ws.Range("E95:E98").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# --- Move the viewport / selection the way the author left it ---
$ws.Range("C55").Select()

Write-Output "done"
